$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.107.40'
$ws.Range("E2").Value = '  -0.70%  '
$ws.Range("D3").Value = '1.912.07'
$ws.Range("E3").Value = '  -1.02%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7403'
$ws.Range("E5").Value = '  -2.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '245.12'
$ws.Range("E6").Value = '  +0.63%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3101'
$ws.Range("E8").Value = '  -2.58%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '26.60'
$ws.Range("E9").Value = '  -4.98%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06988'
$ws.Range("E10").Value = '  -0.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08079'
$ws.Range("E11").Value = '  +0.82%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7729'
$ws.Range("E12").Value = '  -0.85%  '
$ws.Range("D13").Value = '1.895.18'
$ws.Range("E13").Value = '  -1.96%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.351'
$ws.Range("E14").Value = '  -0.72%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.32'
$ws.Range("E15").Value = '  -1.01%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.38'
$ws.Range("E16").Value = '  -0.35%  '
$ws.Range("D17").Value = '30.111.71'
$ws.Range("E17").Value = '  -0.68%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.096'
$ws.Range("E18").Value = '  +2.03%  '
$ws.Range("E19").Value = '  -1.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '240.87'
$ws.Range("E20").Value = '  -4.60%  '
$ws.Range("D21").Value = '2.173.12'
$ws.Range("E21").Value = '  -0.86%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.094'
$ws.Range("E24").Value = '  +5.78%  '
$ws.Range("E25").Value = '  -0.86%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.45'
$ws.Range("E26").Value = '  +1.93%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.00'
$ws.Range("E27").Value = '  -0.34%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1285'
$ws.Range("E28").Value = '  -2.36%  '
$ws.Range("E29").Value = '  -7.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.554'
$ws.Range("E30").Value = '  +2.36%  '
$ws.Range("E31").Value = '  -1.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.349'
$ws.Range("E32").Value = '  -1.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.100'
$ws.Range("E33").Value = '  -0.87%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.313'
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.05152'
$ws.Range("E35").Value = '  -1.40%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7504'
$ws.Range("E36").Value = '  -0.73%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.727'
$ws.Range("E37").Value = '  -2.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01965'
$ws.Range("E38").Value = '  +0.85%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.798'
$ws.Range("E39").Value = '  -0.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.343'
$ws.Range("E40").Value = '  -2.50%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4529'
$ws.Range("E41").Value = '  +0.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '74.73'
$ws.Range("E42").Value = '  -4.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.995'
$ws.Range("E43").Value = '  +1.02%  '
$ws.Range("E44").Value = '  +0.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8392'
$ws.Range("E45").Value = '  +0.27%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.798'
$ws.Range("E46").Value = '  +2.92%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.972'
$ws.Range("E47").Value = '  +0.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '101.92'
$ws.Range("E48").Value = '  +0.49%  '
$ws.Range("D49").Value = '2.060.74'
$ws.Range("E49").Value = '  -1.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.88'
$ws.Range("E50").Value = '  -2.40%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1183'
$ws.Range("E51").Value = '  -2.60%  '
